# Update imputed KNN result values in column E for the result_data_KNN.xlsx workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 16.419
$ws.Range("E12").Value = 17.792
$ws.Range("E18").Value = 16.544
$ws.Range("E37").Value = 16.617
$ws.Range("E55").Value = 16.494
$ws.Range("E68").Value = 17.421
$ws.Range("E77").Value = 16.749
$ws.Range("E78").Value = 16.542
